$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 54
$ws.Range("H54").Value = 20431.1
$ws.Range("I54").Value = 9827.75
$ws.Range("J54").Value = 27500
$ws.Range("K54").Value = 9827.75
$ws.Range("L54").Value = 27500
$ws.Range("M54").Value = -9341.75
$ws.Range("N54").Value = -28472

# Row 99
$ws.Range("H99").Value = 634.4
$ws.Range("I99").Value = 482.66666
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1447.99998
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = 50.00001999999995
$ws.Range("N99").Value = -8996

# Row 100
$ws.Range("H100").Value = 3009.4443
$ws.Range("I100").Value = 2881.4285
$ws.Range("J100").Value = 3090.9092
$ws.Range("K100").Value = 2881.4285
$ws.Range("L100").Value = 3090.9092
$ws.Range("M100").Value = -2340.4285
$ws.Range("N100").Value = -4172.9092

# Row 132
$ws.Range("H132").Value = 756363.5
$ws.Range("I132").Value = 2127.5962
$ws.Range("J132").Value = 3773307.2
$ws.Range("K132").Value = 6382.7886
$ws.Range("L132").Value = 11319921.6
$ws.Range("M132").Value = -3852.7886
$ws.Range("N132").Value = -11324981.6

# Row 135
$ws.Range("H135").Value = 40511.348
$ws.Range("I135").Value = 47169.227
$ws.Range("J135").Value = 3893
$ws.Range("K135").Value = 424523.043
$ws.Range("L135").Value = 35037
$ws.Range("M135").Value = -421988.043
$ws.Range("N135").Value = -40107

# Row 138
$ws.Range("H138").Value = 3849337.2
$ws.Range("I138").Value = 4747
$ws.Range("J138").Value = 4548353.5
$ws.Range("K138").Value = 14241
$ws.Range("L138").Value = 13645060.5
$ws.Range("M138").Value = -9101
$ws.Range("N138").Value = -13655340.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 973.6
$ws.Range("I45").Value = 800.36365
$ws.Range("K45").Value = 800.36365
$ws.Range("M45").Value = -423.36365

# Row 52
$ws.Range("H52").Value = 50780
$ws.Range("J52").Value = 50780
$ws.Range("L52").Value = 50780
$ws.Range("N52").Value = -51416

# Row 61
$ws.Range("H61").Value = 28629826
$ws.Range("I61").Value = 34518230
$ws.Range("J61").Value = 169198.33
$ws.Range("K61").Value = 34518230
$ws.Range("L61").Value = 169198.33
$ws.Range("M61").Value = -34518018
$ws.Range("N61").Value = -169622.33

# Row 102
$ws.Range("H102").Value = 3880
$ws.Range("J102").Value = 4333.3335
$ws.Range("L102").Value = 4333.3335
$ws.Range("N102").Value = -7577.3335

# Row 136
$ws.Range("H136").Value = 28629826
$ws.Range("I136").Value = 34518230
$ws.Range("J136").Value = 169198.33
$ws.Range("K136").Value = 103554690
$ws.Range("L136").Value = 507594.99
$ws.Range("M136").Value = -103552140
$ws.Range("N136").Value = -512694.99

$ws = $wb.Worksheets.Item("BSM")
# Row 140
$ws.Range("H140").Value = 66697.14
$ws.Range("J140").Value = 66697.14
$ws.Range("L140").Value = 66697.14
$ws.Range("N140").Value = -77057.14

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 265068.9
$ws.Range("I31").Value = 47025.863
$ws.Range("J31").Value = 464941.66
$ws.Range("K31").Value = 47025.863
$ws.Range("L31").Value = 464941.66
$ws.Range("M31").Value = -46730.863
$ws.Range("N31").Value = -465531.66

# Row 34
$ws.Range("H34").Value = 265068.9
$ws.Range("I34").Value = 47025.863
$ws.Range("J34").Value = 464941.66
$ws.Range("K34").Value = 47025.863
$ws.Range("L34").Value = 464941.66
$ws.Range("M34").Value = -46823.863
$ws.Range("N34").Value = -465345.66

# Row 86
$ws.Range("H86").Value = 2818.9092
$ws.Range("I86").Value = 2050
$ws.Range("K86").Value = 2050
$ws.Range("M86").Value = -927

# Row 89
$ws.Range("H89").Value = 2818.9092
$ws.Range("I89").Value = 2050
$ws.Range("K89").Value = 10250
$ws.Range("M89").Value = -4634

# Row 99
$ws.Range("H99").Value = 1555.1875
$ws.Range("I99").Value = 1538.3334
$ws.Range("J99").Value = 1576.8572
$ws.Range("K99").Value = 1538.3334
$ws.Range("L99").Value = 1576.8572
$ws.Range("M99").Value = -40.33339999999998
$ws.Range("N99").Value = -4572.8572

# Row 126
$ws.Range("H126").Value = 1555.1875
$ws.Range("I126").Value = 1538.3334
$ws.Range("J126").Value = 1576.8572
$ws.Range("K126").Value = 4615.0002
$ws.Range("L126").Value = 4730.571599999999
$ws.Range("M126").Value = -2145.0002
$ws.Range("N126").Value = -9670.571599999999

# Row 134
$ws.Range("H134").Value = 34648.438
$ws.Range("I134").Value = 753.2857
$ws.Range("J134").Value = 61011.332
$ws.Range("K134").Value = 2259.8571
$ws.Range("L134").Value = 183033.996
$ws.Range("M134").Value = 275.1428999999998
$ws.Range("N134").Value = -188103.996

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 598.03125
$ws.Range("J113").Value = 637.3333
$ws.Range("L113").Value = 1911.9999
$ws.Range("N113").Value = -6251.9999

# Row 122
$ws.Range("H122").Value = 1123.2963
$ws.Range("I122").Value = 398
$ws.Range("J122").Value = 1213.9584
$ws.Range("K122").Value = 3582
$ws.Range("L122").Value = 10925.6256
$ws.Range("M122").Value = -1132
$ws.Range("N122").Value = -15825.6256

$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 50250
$ws.Range("J51").Value = 50250
$ws.Range("L51").Value = 50250
$ws.Range("N51").Value = -51268

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 620.8148
$ws.Range("I22").Value = 529.2308
$ws.Range("J22").Value = 705.8570999999999
$ws.Range("K22").Value = 529.2308
$ws.Range("L22").Value = 705.8570999999999
$ws.Range("M22").Value = -234.2308
$ws.Range("N22").Value = -1295.8571

# Row 27
$ws.Range("H27").Value = 620.8148
$ws.Range("I27").Value = 529.2308
$ws.Range("J27").Value = 705.8570999999999
$ws.Range("K27").Value = 529.2308
$ws.Range("L27").Value = 705.8570999999999
$ws.Range("M27").Value = -422.2308
$ws.Range("N27").Value = -919.8570999999999

# Row 61
$ws.Range("H61").Value = 2150.4211
$ws.Range("I61").Value = 2129.8667
$ws.Range("J61").Value = 2227.5
$ws.Range("K61").Value = 2129.8667
$ws.Range("L61").Value = 2227.5
$ws.Range("M61").Value = -1927.8667
$ws.Range("N61").Value = -2631.5

# Row 82
$ws.Range("H82").Value = 1193.5
$ws.Range("I82").Value = 1193.5
$ws.Range("K82").Value = 1193.5
$ws.Range("M82").Value = -832.5

# Row 85
$ws.Range("H85").Value = 1193.5
$ws.Range("I85").Value = 1193.5
$ws.Range("K85").Value = 1193.5
$ws.Range("M85").Value = 54.5

# Row 93
$ws.Range("H93").Value = 1451.9412
$ws.Range("I93").Value = 1405.9286
$ws.Range("J93").Value = 1666.6666
$ws.Range("K93").Value = 1405.9286
$ws.Range("L93").Value = 1666.6666
$ws.Range("M93").Value = -157.9286
$ws.Range("N93").Value = -4162.6666

# Row 100
$ws.Range("H100").Value = 1927.5454
$ws.Range("I100").Value = 1840.6
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1840.6
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1299.6
$ws.Range("N100").Value = -3082

# Row 113
$ws.Range("H113").Value = 2150.4211
$ws.Range("I113").Value = 2129.8667
$ws.Range("J113").Value = 2227.5
$ws.Range("K113").Value = 2129.8667
$ws.Range("L113").Value = 2227.5
$ws.Range("M113").Value = 40.13329999999996
$ws.Range("N113").Value = -6567.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 389.42105
$ws.Range("I107").Value = 323.96667
$ws.Range("K107").Value = 971.9000100000001
$ws.Range("M107").Value = 948.0999899999999

# Row 113
$ws.Range("H113").Value = 544.46344
$ws.Range("I113").Value = 656.5417
$ws.Range("J113").Value = 386.2353
$ws.Range("K113").Value = 1969.6251
$ws.Range("L113").Value = 1158.7059
$ws.Range("M113").Value = 200.3749
$ws.Range("N113").Value = -5498.7059

Write-Host "Applied all changes"
